$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data region (rows 2-63) before rewriting with the corrected dataset
$ws.Range("A2:C63").ClearContents()

$ws.Range("A2").Value = "-"
$ws.Range("B2").Value = "section"
$ws.Range("C2").Value = "Structure Preparation"

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "examined organisms"
$ws.Range("C3").Value = "mouse and rat"

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "stage"
$ws.Range("C4").Value = "sequence alignment"

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "sequence alignment"
$ws.Range("C5").Value = "GluN1,GluN2,GluN3"

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Software"
$ws.Range("C6").Value = "MAAFT Server"

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "version"
$ws.Range("C7").Value = "'7"

$ws.Range("A8").Value = 2
$ws.Range("B8").Value = "settings"
$ws.Range("C8").Value = "default"

$ws.Range("A9").Value = 3
$ws.Range("B9").Value = "stage"
$ws.Range("C9").Value = "alignment refinement"

$ws.Range("A10").Value = 3
$ws.Range("B10").Value = "software"
$ws.Range("C10").Value = "GLProbs"

$ws.Range("A11").Value = 3
$ws.Range("B11").Value = "consistency transformation passes"
$ws.Range("C11").Value = "'2"

$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "iterative refinement passes"
$ws.Range("C12").Value = "'100"

$ws.Range("A13").Value = 4
$ws.Range("B13").Value = "stage"
$ws.Range("C13").Value = "manual sequence alignment"

$ws.Range("A14").Value = 4
$ws.Range("B14").Value = "base template"
$ws.Range("C14").Value = "result of previous step"

$ws.Range("A15").Value = 4
$ws.Range("B15").Value = "template structure"
$ws.Range("C15").Value = "GluN2DLBD"

$ws.Range("A16").Value = 4
$ws.Range("B16").Value = "PDB ID"
$ws.Range("C16").Value = "3OEK"

$ws.Range("A17").Value = 4
$ws.Range("B17").Value = "template structure"
$ws.Range("C17").Value = "GluN3BLBD"

$ws.Range("A18").Value = 4
$ws.Range("B18").Value = "PDB ID"
$ws.Range("C18").Value = "2RCA"

$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "software"
$ws.Range("C19").Value = "BioLuminate package"

$ws.Range("A20").Value = 5
$ws.Range("B20").Value = "step type"
$ws.Range("C20").Value = "iteration"

$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "flow type"
$ws.Range("C21").Value = "for each"

$ws.Range("A22").Value = 5
$ws.Range("B22").Value = "flow parameter"
$ws.Range("C22").Value = " generated pose"

$ws.Range("A23").Value = 5
$ws.Range("B23").Value = "operation"
$ws.Range("C23").Value = "energy minimization"

$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "(minimization) target"
$ws.Range("C24").Value = "receptor residue"

$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "target criteria"
$ws.Range("C25").Value = "lte 5  Å"

$ws.Range("A26").Value = 5
$ws.Range("B26").Value = "operation"
$ws.Range("C26").Value = "optimization"

$ws.Range("A27").Value = 5
$ws.Range("B27").Value = "(optimization) target"
$ws.Range("C27").Value = "side chain rotamers"

$ws.Range("A28").Value = 6
$ws.Range("B28").Value = "step type"
$ws.Range("C28").Value = "iteration"

$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "flow type"
$ws.Range("C29").Value = "while"

$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "flow parameter"
$ws.Range("C30").Value = "pH"

$ws.Range("A31").Value = 6
$ws.Range("B31").Value = "flow logical parameter"
$ws.Range("C31").Value = "lte"

$ws.Range("A32").Value = 6
$ws.Range("B32").Value = "flow compared value"
$ws.Range("C32").Value = "'7"

$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "flow operation"
$ws.Range("C33").Value = "+"

$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "flow magnitude"
$ws.Range("C34").Value = "'1"

$ws.Range("A35").Value = 6
$ws.Range("B35").Value = "operation"
$ws.Range("C35").Value = "energy minimization"

$ws.Range("A36").Value = 6
$ws.Range("B36").Value = "(minimization)  target"
$ws.Range("C36").Value = "receptor residue"

$ws.Range("A37").Value = 6
$ws.Range("B37").Value = "target criteria"
$ws.Range("C37").Value = "lte 5  Å"

$ws.Range("A38").Value = 6
$ws.Range("B38").Value = "operation"
$ws.Range("C38").Value = "optimization"

$ws.Range("A39").Value = 6
$ws.Range("B39").Value = "(optimization) target"
$ws.Range("C39").Value = "side chain rotamers"

$ws.Range("A40").Value = 7
$ws.Range("B40").Value = "step type"
$ws.Range("C40").Value = "conditional"

$ws.Range("A41").Value = 7
$ws.Range("B41").Value = "flow type"
$ws.Range("C41").Value = "if"

$ws.Range("A42").Value = 7
$ws.Range("B42").Value = "flow parameter"
$ws.Range("C42").Value = "pH"

$ws.Range("A43").Value = 7
$ws.Range("B43").Value = "flow logical parameter"
$ws.Range("C43").Value = "lte"

$ws.Range("A44").Value = 7
$ws.Range("B44").Value = "flow compared value"
$ws.Range("C44").Value = "'7 "

$ws.Range("A45").Value = 7
$ws.Range("B45").Value = "operation"
$ws.Range("C45").Value = "energy minimization"

$ws.Range("A46").Value = 7
$ws.Range("B46").Value = "(minimization) target"
$ws.Range("C46").Value = "receptor residue"

$ws.Range("A47").Value = 7
$ws.Range("B47").Value = "target criteria"
$ws.Range("C47").Value = "lte 5  Å"

$ws.Range("A48").Value = 7
$ws.Range("B48").Value = "operation"
$ws.Range("C48").Value = "optimization"

$ws.Range("A49").Value = 7
$ws.Range("B49").Value = "(optimization) target"
$ws.Range("C49").Value = "side chain rotamers"

$ws.Range("A50").Value = 8
$ws.Range("B50").Value = "step type"
$ws.Range("C50").Value = "conditional"

$ws.Range("A51").Value = 8
$ws.Range("B51").Value = "flow type"
$ws.Range("C51").Value = "else if"

$ws.Range("A52").Value = 8
$ws.Range("B52").Value = "flow parameter"
$ws.Range("C52").Value = "pH"

$ws.Range("A53").Value = 8
$ws.Range("B53").Value = "flow logical parameter"
$ws.Range("C53").Value = "between"

$ws.Range("A54").Value = 8
$ws.Range("B54").Value = "flow range"
$ws.Range("C54").Value = "[8-12]"

$ws.Range("A55").Value = 8
$ws.Range("B55").Value = "start iteration value"
$ws.Range("C55").Value = 8

$ws.Range("A56").Value = 8
$ws.Range("B56").Value = "end iteration value"
$ws.Range("C56").Value = 12

$ws.Range("A57").Value = 8
$ws.Range("B57").Value = "operation"
$ws.Range("C57").Value = "energy minimization"

$ws.Range("A58").Value = 8
$ws.Range("B58").Value = "(minimization) target"
$ws.Range("C58").Value = "receptor residue"

$ws.Range("A59").Value = 8
$ws.Range("B59").Value = "target criteria"
$ws.Range("C59").Value = "lte 5  Å"

$ws.Range("A60").Value = 8
$ws.Range("B60").Value = "operation"
$ws.Range("C60").Value = "optimization"

$ws.Range("A61").Value = 8
$ws.Range("B61").Value = "(optimization) target"
$ws.Range("C61").Value = "side chain rotamers"

$ws.Range("A62").Value = 9
$ws.Range("B62").Value = "step type"
$ws.Range("C62").Value = "conditional"

$ws.Range("A63").Value = 9
$ws.Range("B63").Value = "flow type"
$ws.Range("C63").Value = "else"

$ws.Range("A64").Value = 9
$ws.Range("B64").Value = "operation"
$ws.Range("C64").Value = "energy minimization"

$ws.Range("A65").Value = 9
$ws.Range("B65").Value = "(minimization) target"
$ws.Range("C65").Value = "receptor residue"

$ws.Range("A66").Value = 9
$ws.Range("B66").Value = "target criteria"
$ws.Range("C66").Value = "lte 5  Å"

$ws.Range("A67").Value = 9
$ws.Range("B67").Value = "operation"
$ws.Range("C67").Value = "optimization"

$ws.Range("A68").Value = 9
$ws.Range("B68").Value = "(optimization) target"
$ws.Range("C68").Value = "side chain rotamers"

$ws.Range("A69").Value = 10
$ws.Range("B69").Value = "step type"
$ws.Range("C69").Value = "iteration"

$ws.Range("A70").Value = 10
$ws.Range("B70").Value = "flow type"
$ws.Range("C70").Value = "for"

$ws.Range("A71").Value = 10
$ws.Range("B71").Value = "flow parameter"
$ws.Range("C71").Value = "pH"

$ws.Range("A72").Value = 10
$ws.Range("B72").Value = "flow range"
$ws.Range("C72").Value = "[1-7]"

$ws.Range("A73").Value = 10
$ws.Range("B73").Value = "start iteration value"
$ws.Range("C73").Value = 1

$ws.Range("A74").Value = 10
$ws.Range("B74").Value = "end iteration value"
$ws.Range("C74").Value = 7

$ws.Range("A75").Value = 10
$ws.Range("B75").Value = "flow operation"
$ws.Range("C75").Value = "+"

$ws.Range("A76").Value = 10
$ws.Range("B76").Value = "flow magnitude"
$ws.Range("C76").Value = "'1"

$ws.Range("A77").Value = 10
$ws.Range("B77").Value = "operation"
$ws.Range("C77").Value = "energy minimization"

$ws.Range("A78").Value = 10
$ws.Range("B78").Value = "(minimization) target"
$ws.Range("C78").Value = "receptor residue"

$ws.Range("A79").Value = 10
$ws.Range("B79").Value = "target criteria"
$ws.Range("C79").Value = "lte 5  Å"

$ws.Range("A80").Value = 10
$ws.Range("B80").Value = "operation"
$ws.Range("C80").Value = "optimization"

$ws.Range("A81").Value = 10
$ws.Range("B81").Value = "(optimization) target"
$ws.Range("C81").Value = "side chain rotamers"
